# Update "想去人数" (interest count) values in column F for rows 4-6
# on both the "展览" sheet and the "全部类型" sheet.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 103
    $ws.Range("F5").Value = 2669
    $ws.Range("F6").Value = 252
}
